$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("FEINmismatch")
$ws2 = $wb.Worksheets.Item("FEINSSNmismatch")

# New timestamp strings generated by the latest RAD test run (Jan 28 2025),
# replacing the previous run's timestamps (Nov 12 2024) in column B.
$sheet1Timestamps = @(
    "Tue Jan 28 21:27:27 EST 2025",
    "Tue Jan 28 21:27:42 EST 2025",
    "Tue Jan 28 21:27:54 EST 2025",
    "Tue Jan 28 21:28:06 EST 2025",
    "Tue Jan 28 21:28:19 EST 2025",
    "Tue Jan 28 21:28:31 EST 2025",
    "Tue Jan 28 21:28:42 EST 2025",
    "Tue Jan 28 21:28:53 EST 2025",
    "Tue Jan 28 21:29:04 EST 2025",
    "Tue Jan 28 21:29:15 EST 2025",
    "Tue Jan 28 21:29:26 EST 2025",
    "Tue Jan 28 21:29:37 EST 2025",
    "Tue Jan 28 21:29:48 EST 2025",
    "Tue Jan 28 21:29:59 EST 2025",
    "Tue Jan 28 21:30:11 EST 2025",
    "Tue Jan 28 21:30:22 EST 2025",
    "Tue Jan 28 21:30:33 EST 2025",
    "Tue Jan 28 21:30:44 EST 2025",
    "Tue Jan 28 21:30:55 EST 2025",
    "Tue Jan 28 21:31:06 EST 2025",
    "Tue Jan 28 21:31:17 EST 2025",
    "Tue Jan 28 21:31:28 EST 2025",
    "Tue Jan 28 21:31:39 EST 2025",
    "Tue Jan 28 21:31:50 EST 2025",
    "Tue Jan 28 21:32:01 EST 2025",
    "Tue Jan 28 21:32:12 EST 2025",
    "Tue Jan 28 21:32:23 EST 2025",
    "Tue Jan 28 21:32:34 EST 2025",
    "Tue Jan 28 21:32:45 EST 2025"
)

$sheet2Timestamps = @(
    "Tue Jan 28 21:32:56 EST 2025",
    "Tue Jan 28 21:33:07 EST 2025",
    "Tue Jan 28 21:33:18 EST 2025",
    "Tue Jan 28 21:33:29 EST 2025",
    "Tue Jan 28 21:33:39 EST 2025",
    "Tue Jan 28 21:33:50 EST 2025",
    "Tue Jan 28 21:34:00 EST 2025",
    "Tue Jan 28 21:34:11 EST 2025",
    "Tue Jan 28 21:34:22 EST 2025",
    "Tue Jan 28 21:34:32 EST 2025",
    "Tue Jan 28 21:34:43 EST 2025",
    "Tue Jan 28 21:34:53 EST 2025",
    "Tue Jan 28 21:35:04 EST 2025",
    "Tue Jan 28 21:35:15 EST 2025",
    "Tue Jan 28 21:35:25 EST 2025",
    "Tue Jan 28 21:35:36 EST 2025",
    "Tue Jan 28 21:35:47 EST 2025",
    "Tue Jan 28 21:35:57 EST 2025"
)

for ($i = 0; $i -lt $sheet1Timestamps.Length; $i++) {
    $row = $i + 2
    $ws1.Range("B$row").Value = $sheet1Timestamps[$i]
}

for ($i = 0; $i -lt $sheet2Timestamps.Length; $i++) {
    $row = $i + 2
    $ws2.Range("B$row").Value = $sheet2Timestamps[$i]
}

# Rows whose Result flipped from Pass to Fail in the new run.
$ws1.Range("A18").Value = "Fail"
$ws1.Range("A19").Value = "Fail"
$ws1.Range("A29").Value = "Fail"
